# Apply updated "K" (column G) values to the save_data sheet.
# The K column values were regenerated (previously "Strike#"-derived values),
# so only the numeric contents of G2:G15 change; everything else is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 3
    3  = 3
    4  = 2
    5  = 4
    6  = 5
    7  = 2
    8  = 3
    9  = 1
    10 = 3
    11 = 2
    12 = 3
    13 = 2
    14 = 4
    15 = 4
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
